$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Mark the three inline pictures as "no proofing" (<w:noProof/> in rPr).
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $shp = $d.InlineShapes.Item($i)
    $shp.Range.NoProofing = 1
}

# ---------------------------------------------------------------------------
# 2) Fix the duplicated word "menang menang" -> "menang" in the Atalanta /
#    AC Milan question, while splitting the sentence into three runs
#    (mirrors the interactive edit-in-place behaviour of Word).
# ---------------------------------------------------------------------------
$fullText = $d.Content.Text
$beforeDup = "Untuk setiap kemenangan, klub akan mendapatkan nilai 3 poin, imbang 1 poin dan kalah 0 poin. Jika pertandingan tersisa Atalanta memenangkan dua pertandingan dan sisanya imbang, kemungkinan komposisi menang"
$idx = $fullText.IndexOf($beforeDup)
if ($idx -ge 0) {
    $cutStart = $idx + $beforeDup.Length
    $cutEnd = $cutStart + 8   # deletes " menang " (space + "menang" + space)

    $rngDel = $d.Range($cutStart, $cutEnd)
    $rngDel.Text = ""

    # Re-insert a single space as its own run, nudging formatting so the
    # engine keeps it as a distinct run instead of re-merging it with its
    # neighbours.
    $spaceRange = $d.Range($cutStart, $cutStart)
    $spaceRange.InsertAfter(" ")
    $spaceRange2 = $d.Range($cutStart, $cutStart + 1)
    $spaceRange2.Bold = 1
    $spaceRange2.Bold = 0
}

# ---------------------------------------------------------------------------
# 3) Fill in the answer key table (previously-empty answer cells).
# ---------------------------------------------------------------------------
$t = $d.Tables.Item(1)

function Set-AnswerCell($row, $col, $text) {
    $cell = $t.Cell($row, $col)
    $cell.Range.Text = $text
    $cell.Range.Font.Name = "Times New Roman"
    $cell.Range.Font.Size = 12
}

Set-AnswerCell 1 4 "C"
Set-AnswerCell 2 4 "A"
Set-AnswerCell 3 2 "A"
Set-AnswerCell 4 2 "B"
Set-AnswerCell 5 2 "D"

Write-Output "edit complete"
